# Generate Report for Handback
#
# This applies the "handback" pass of the localization-status report:
#  - the overall Status text flips from "Ready for handoff" to
#    "Handed back: in sync with en-US" everywhere it is used,
#  - the per-language sheets (zh-cn, de-de) get their "Latest Target File"
#    and "Latest Handback File" columns populated (with a hyperlink on the
#    target-file cell, mirroring the existing Source-File-Name hyperlink),
#  - the "Latest Handback DateTime" stamps are written,
#  - a handful of columns are widened to fit the newly-populated content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (every cell that shares this string has to be rewritten individually)
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# 2. Column widths - widen to fit the now-populated / longer content.
#    (ColumnWidth is expressed in characters; Excel snaps it to whole
#    pixels, so these are the inputs that land closest to the target.)
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668   # E
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668   # F

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.166666666666668      # C
$wsZhCn.Columns.Item(9).ColumnWidth  = 39.166666666666664      # I
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664      # J

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.166666666666668      # C
$wsDeDe.Columns.Item(9).ColumnWidth  = 39.166666666666664      # I
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664      # J

# ---------------------------------------------------------------------
# 3. Helper: populate "Latest Target File" (I) with a hyperlink to the
#    source .md file (just like column A already has), "Latest Handback
#    File" (J) with the generated xlf name, and "Latest Handback
#    DateTime" (K) with the handback timestamp.
# ---------------------------------------------------------------------
function Set-HandbackRow {
    param($ws, $row, $mdName, $xlfName, $dateTime, $rId)

    $url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/45bb93908cbba444cceaab03faf865e08e1e4f62/e2e/$mdName"

    $targetCell = $ws.Cells.Item($row, 9)   # column I
    $ws.Hyperlinks.Add($targetCell, $url, "", "", $mdName)
    $targetCell.Font.Underline = 2
    $targetCell.Font.Color = 15570276

    $ws.Cells.Item($row, 10).Value = $xlfName   # column J
    $ws.Cells.Item($row, 11).Value = $dateTime  # column K
}

# zh-cn
Set-HandbackRow $wsZhCn 2 "32c5c11f-5cd8-4ed8-b45a-5684172a7f1b.md" `
    "32c5c11f-5cd8-4ed8-b45a-5684172a7f1b.d7d73d8018e147d6174f6e782aaa0d9d783b2ea9.zh-cn.xlf" `
    "2016-08-26 10:46:06"

Set-HandbackRow $wsZhCn 3 "63679998-18cd-485d-97e6-76d9b68749bd.md" `
    "63679998-18cd-485d-97e6-76d9b68749bd.134fa91a48401a5a341eeb3756855d81b3abc47e.zh-cn.xlf" `
    "2016-08-26 10:46:06"

# de-de
Set-HandbackRow $wsDeDe 2 "32c5c11f-5cd8-4ed8-b45a-5684172a7f1b.md" `
    "32c5c11f-5cd8-4ed8-b45a-5684172a7f1b.d7d73d8018e147d6174f6e782aaa0d9d783b2ea9.de-de.xlf" `
    "2016-08-26 10:46:15"

Set-HandbackRow $wsDeDe 3 "63679998-18cd-485d-97e6-76d9b68749bd.md" `
    "63679998-18cd-485d-97e6-76d9b68749bd.134fa91a48401a5a341eeb3756855d81b3abc47e.de-de.xlf" `
    "2016-08-26 10:46:15"
